# Formulas/score sheet: mark additional "Раздел" rows with a pass/fail
# flag in column N — "+" for included articles, "-" for excluded ones.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист1")

$ws.Range("N3").Value = "+"
$ws.Range("N4").Value = "+"
$ws.Range("N6").Value = "-"
$ws.Range("N9").Value = "-"
$ws.Range("N12").Value = "+"
$ws.Range("N17").Value = "-"
$ws.Range("N18").Value = "+"

# Move the active selection to N5, matching where the author left off editing.
[void]$ws.Range("N5").Select()
